$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 487.5
$ws.Range("J7").Value = 487.5
$ws.Range("L7").Value = 487.5
$ws.Range("N7").Value = -711.5

$ws.Range("H14").Value = 487.5
$ws.Range("J14").Value = 487.5
$ws.Range("L14").Value = 487.5
$ws.Range("N14").Value = -869.5

$ws.Range("H40").Value = 1758
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1758
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1758
$ws.Range("N40").Value = -2108
$ws.Range("M40").Value = $null

$ws.Range("H74").Value = 2609.1667
$ws.Range("I74").Value = 2468.5293
$ws.Range("K74").Value = 2468.5293
$ws.Range("M74").Value = -1532.5293

$ws.Range("H77").Value = 2609.1667
$ws.Range("I77").Value = 2468.5293
$ws.Range("K77").Value = 12342.6465
$ws.Range("M77").Value = -7662.646500000001

$ws.Range("H100").Value = 3999.6
$ws.Range("I100").Value = 3999.6
$ws.Range("K100").Value = 3999.6
$ws.Range("M100").Value = -3458.6

$ws.Range("H107").Value = 4443.3335
$ws.Range("I107").Value = 1665
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 1665
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = 255
$ws.Range("N107").Value = -13840

$ws.Range("H112").Value = 3274.75
$ws.Range("J112").Value = 3966.3333
$ws.Range("L112").Value = 11898.9999
$ws.Range("N112").Value = -14114.9999

$ws.Range("H132").Value = 2344.2593
$ws.Range("I132").Value = 2344.2593
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7032.777900000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4502.777900000001
$ws.Range("N132").Value = $null

$ws.Range("H141").Value = 3215.6667
$ws.Range("I141").Value = 3215.6667
$ws.Range("K141").Value = 9647.000100000001
$ws.Range("M141").Value = -4467.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2805.6924
$ws.Range("J88").Value = 3555
$ws.Range("L88").Value = 3555
$ws.Range("N88").Value = -4367

$ws.Range("H91").Value = 2805.6924
$ws.Range("J91").Value = 3555
$ws.Range("L91").Value = 3555
$ws.Range("N91").Value = -6363

$ws.Range("H95").Value = 10551.25
$ws.Range("I95").Value = 15000
$ws.Range("J95").Value = 9068.333000000001
$ws.Range("K95").Value = 15000
$ws.Range("L95").Value = 9068.333000000001
$ws.Range("M95").Value = -12254
$ws.Range("N95").Value = -14560.333

$ws.Range("H97").Value = 916.6923
$ws.Range("I97").Value = 918.9167
$ws.Range("J97").Value = 890
$ws.Range("K97").Value = 918.9167
$ws.Range("L97").Value = 890
$ws.Range("M97").Value = -422.9167
$ws.Range("N97").Value = -1882

$ws.Range("H98").Value = 39997.5
$ws.Range("J98").Value = 39997.5
$ws.Range("L98").Value = 39997.5
$ws.Range("N98").Value = -45987.5

$ws.Range("H102").Value = 1996.1666
$ws.Range("I102").Value = 1269
$ws.Range("K102").Value = 1269
$ws.Range("M102").Value = 353

$ws.Range("H104").Value = 70000
$ws.Range("J104").Value = 70000
$ws.Range("L104").Value = 70000
$ws.Range("N104").Value = -76988

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 10000
$ws.Range("K26").Value = 10000
$ws.Range("M26").Value = -9708

$ws.Range("H94").Value = 1559.1852
$ws.Range("I94").Value = 1687.4584
$ws.Range("J94").Value = 533
$ws.Range("K94").Value = 1687.4584
$ws.Range("L94").Value = 533
$ws.Range("M94").Value = -1236.4584
$ws.Range("N94").Value = -1435

$ws.Range("H96").Value = 5428
$ws.Range("I96").Value = 5428
$ws.Range("K96").Value = 5428
$ws.Range("M96").Value = -2682

$ws.Range("H105").Value = 1828.1538
$ws.Range("I105").Value = 1498.6
$ws.Range("K105").Value = 1498.6
$ws.Range("M105").Value = 248.4000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2258.6382
$ws.Range("I31").Value = 1793.1666
$ws.Range("K31").Value = 1793.1666
$ws.Range("M31").Value = -1498.1666

$ws.Range("H34").Value = 2258.6382
$ws.Range("I34").Value = 1793.1666
$ws.Range("K34").Value = 1793.1666
$ws.Range("M34").Value = -1591.1666

$ws.Range("H62").Value = 4112.7144
$ws.Range("I62").Value = 4498
$ws.Range("J62").Value = 3823.75
$ws.Range("K62").Value = 4498
$ws.Range("L62").Value = 3823.75
$ws.Range("M62").Value = -3874
$ws.Range("N62").Value = -5071.75

$ws.Range("H65").Value = 4112.7144
$ws.Range("I65").Value = 4498
$ws.Range("J65").Value = 3823.75
$ws.Range("K65").Value = 22490
$ws.Range("L65").Value = 19118.75
$ws.Range("M65").Value = -19370
$ws.Range("N65").Value = -25358.75

$ws.Range("H105").Value = 2000
$ws.Range("J105").Value = 1500
$ws.Range("L105").Value = 1500
$ws.Range("N105").Value = -4994

$ws.Range("H132").Value = 1441.8572
$ws.Range("I132").Value = 1441.8572
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4325.571599999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1795.571599999999
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 463.84616
$ws.Range("I2").Value = 325.8
$ws.Range("J2").Value = 550.125
$ws.Range("K2").Value = 1954.8
$ws.Range("L2").Value = 3300.75
$ws.Range("M2").Value = -1841.8
$ws.Range("N2").Value = -3526.75

$ws.Range("H7").Value = 900

$ws.Range("H34").Value = 1235.5714
$ws.Range("I34").Value = 616.6667
$ws.Range("J34").Value = 1699.75
$ws.Range("K34").Value = 1850.0001
$ws.Range("L34").Value = 5099.25
$ws.Range("M34").Value = -1766.0001
$ws.Range("N34").Value = -5267.25

$ws.Range("H68").Value = 1302.1428
$ws.Range("J68").Value = 959.6667
$ws.Range("L68").Value = 2879.0001
$ws.Range("N68").Value = -4501.0001

$ws.Range("H71").Value = 1302.1428
$ws.Range("J71").Value = 959.6667
$ws.Range("L71").Value = 8637.0003
$ws.Range("N71").Value = -16749.0003

$ws.Range("H129").Value = 1168.75
$ws.Range("J129").Value = 5997
$ws.Range("L129").Value = 17991
$ws.Range("N129").Value = -27991

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4246.8887
$ws.Range("I102").Value = 4215.25
$ws.Range("K102").Value = 4215.25
$ws.Range("M102").Value = -2593.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2763.7144
$ws.Range("I7").Value = 1874
$ws.Range("K7").Value = 1874
$ws.Range("M7").Value = -1762

$ws.Range("H13").Value = 3402.8
$ws.Range("I13").Value = 3402.8
$ws.Range("K13").Value = 3402.8
$ws.Range("M13").Value = -3262.8

$ws.Range("H93").Value = 1200
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 1200
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 1200
$ws.Range("N93").Value = -3696
$ws.Range("M93").Value = $null

$ws.Range("H108").Value = 59997
$ws.Range("J108").Value = 59997
$ws.Range("L108").Value = 59997
$ws.Range("N108").Value = -67677

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = $null

$ws.Range("H123").Value = 78996.5
$ws.Range("J123").Value = 78996.5
$ws.Range("L123").Value = 78996.5
$ws.Range("N123").Value = -88796.5

$ws.Range("H126").Value = 2763.7144
$ws.Range("I126").Value = 1874
$ws.Range("K126").Value = 5622
$ws.Range("M126").Value = -3152

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4949
$ws.Range("I62").Value = 4898
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 4898
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -4274
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 4949
$ws.Range("I65").Value = 4898
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 24490
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -21370
$ws.Range("N65").Value = -31240

$ws.Range("H96").Value = 300
$ws.Range("I96").Value = 300
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 300
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 1073
$ws.Range("N96").Value = $null

$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140
